$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 header row shrinks back down (author trimmed the wrapped header height)
$ws.Rows.Item(5).RowHeight = 60

# Row 6 - Regular US Data
$ws.Range("B6").Value = 8
$ws.Range("C6").Value = 12
$ws.Range("D6").Value = 126
$ws.Range("E6").Value = 333100360
$ws.Range("F6").Value = 5
$ws.Range("G6").Formula = "=((`$B`$2/8)-(`$B`$2/12)+(`$B`$2/126))*5"
$ws.Range("H6").Formula = "=`$E`$6+`$G`$6"
$ws.Range("I6").Formula = "=IF(`$H6>`$E6, ""Increase"", ""Decrease"")"

# Row 7 - Population Increase - High Birth Rate
$ws.Range("B7").Value = 2
$ws.Range("C7").Value = 30
$ws.Range("D7").Value = 300
$ws.Range("E7").Value = 333100360
$ws.Range("F7").Value = 5
$ws.Range("G7").Formula = "=((`$B2/`$B7)+(`$B2/`$D7)-(`$B2/`$C7))*5"
$ws.Range("H7").Formula = "=`$E7+`$G7"
$ws.Range("I7").Formula = "=IF(`$H7>`$E7, ""Increase"", ""Decrease"")"

# Row 8 - Population Increase - High Migration
$ws.Range("B8").Value = 50
$ws.Range("C8").Value = 50
$ws.Range("D8").Value = 40
$ws.Range("E8").Value = 333100360
$ws.Range("F8").Value = 5
$ws.Range("G8").Formula = "=((`$B2/`$B8)+(`$B2/`$D8)-(`$B2/`$C8))*5"
$ws.Range("H8").Formula = "=`$G8+`$E8"
$ws.Range("I8").Formula = "=IF(`$H8>`$E8, ""Increase"", ""Decrease"")"

# Row 9 - Population Decrease - High Death Rate
$ws.Range("B9").Value = 100
$ws.Range("C9").Value = 3
$ws.Range("D9").Value = 800
$ws.Range("E9").Value = 333100360
$ws.Range("F9").Value = 5
$ws.Range("G9").Formula = "=((`$B2/`$B9)+(`$B2/`$D9)-(`$B2/`$C9))*5"
$ws.Range("H9").Formula = "=E9+G9"
$ws.Range("I9").Formula = "=IF(H9>E9, ""Increase"", ""Decrease"")"

# Row 10 - Population Low Birth Rate and Low Migration
$ws.Range("B10").Value = 70
$ws.Range("C10").Value = 12
$ws.Range("D10").Value = 400
$ws.Range("E10").Value = 333100360
$ws.Range("F10").Value = 5
$ws.Range("G10").Formula = "=((B2/B10)+(B2/D10)-(B2/C10))*5"
$ws.Range("H10").Formula = "=E10+G10"
$ws.Range("I10").Formula = "=IF(H10>E10, ""Increase"", ""Decrease"")"

# Selection moved off the filled table after data entry
$ws.Range("G18").Select() | Out-Null
